$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 20318
$ws.Range("I92").Value = 25272.5
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 25272.5
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = -24024.5
$ws.Range("N92").Value = -2996

$ws.Range("H136").Value = 60580
$ws.Range("J136").Value = 60580
$ws.Range("L136").Value = 60580
$ws.Range("N136").Value = -70780

$ws.Range("H138").Value = 2024976.5
$ws.Range("J138").Value = 6468.904
$ws.Range("L138").Value = 19406.712
$ws.Range("N138").Value = -29686.712

$ws.Range("H139").Value = 72635
$ws.Range("J139").Value = 72635
$ws.Range("L139").Value = 72635
$ws.Range("N139").Value = -82915

$ws.Range("H140").Value = 111275
$ws.Range("I140").Value = 93550
$ws.Range("J140").Value = 114820
$ws.Range("K140").Value = 93550
$ws.Range("L140").Value = 114820
$ws.Range("M140").Value = -88370
$ws.Range("N140").Value = -125180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 112644.445
$ws.Range("I2").Value = 1200
$ws.Range("J2").Value = 144485.72
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 144485.72
$ws.Range("M2").Value = -1087
$ws.Range("N2").Value = -144711.72

$ws.Range("H53").Value = 14034.5
$ws.Range("J53").Value = 27600
$ws.Range("L53").Value = 27600
$ws.Range("N53").Value = -28964

$ws.Range("H110").Value = 895.1177
$ws.Range("I110").Value = 821.13336
$ws.Range("J110").Value = 1450
$ws.Range("K110").Value = 821.13336
$ws.Range("L110").Value = 1450
$ws.Range("M110").Value = 1223.86664
$ws.Range("N110").Value = -5540

$ws.Range("H116").Value = 112644.445
$ws.Range("I116").Value = 1200
$ws.Range("J116").Value = 144485.72
$ws.Range("K116").Value = 1200
$ws.Range("L116").Value = 144485.72
$ws.Range("M116").Value = 1094
$ws.Range("N116").Value = -149073.72

$ws.Range("H133").Value = 31638.092
$ws.Range("J133").Value = 31638.092
$ws.Range("L133").Value = 31638.092
$ws.Range("N133").Value = -36698.092

$ws.Range("H134").Value = 53181.668
$ws.Range("J134").Value = 53181.668
$ws.Range("L134").Value = 53181.668
$ws.Range("N134").Value = -63321.668

$ws.Range("H135").Value = 48659.89
$ws.Range("I135").Value = 47390
$ws.Range("J135").Value = 48818.625
$ws.Range("K135").Value = 47390
$ws.Range("L135").Value = 48818.625
$ws.Range("M135").Value = -42320
$ws.Range("N135").Value = -58958.625

$ws.Range("H138").Value = 62600
$ws.Range("J138").Value = 62600
$ws.Range("L138").Value = 62600
$ws.Range("N138").Value = -72880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 112644.445
$ws.Range("I3").Value = 1200
$ws.Range("J3").Value = 144485.72
$ws.Range("K3").Value = 1200
$ws.Range("L3").Value = 144485.72
$ws.Range("M3").Value = -1086
$ws.Range("N3").Value = -144713.72

$ws.Range("H105").Value = 3619.2222
$ws.Range("I105").Value = 3619.2222
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3619.2222
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1872.2222
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 23000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 32000
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 32000
$ws.Range("M51").Value = -4264
$ws.Range("N51").Value = -33472

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 23000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 32000
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 32000
$ws.Range("M61").Value = -4652
$ws.Range("N61").Value = -32696

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 44999
$ws.Range("J131").Value = 44999
$ws.Range("L131").Value = 44999
$ws.Range("N131").Value = -55079

$ws.Range("H137").Value = 74780
$ws.Range("J137").Value = 74780
$ws.Range("L137").Value = 74780
$ws.Range("N137").Value = -84980

$ws.Range("H138").Value = 73641.42999999999
$ws.Range("J138").Value = 73641.42999999999
$ws.Range("L138").Value = 73641.42999999999
$ws.Range("N138").Value = -83921.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 76862.86
$ws.Range("I4").Value = 91007.27
$ws.Range("J4").Value = 25000
$ws.Range("K4").Value = 273021.81
$ws.Range("L4").Value = 75000
$ws.Range("M4").Value = -272909.81
$ws.Range("N4").Value = -75224

$ws.Range("H56").Value = 5866.1665
$ws.Range("I56").Value = 5866.1665
$ws.Range("K56").Value = 5866.1665
$ws.Range("M56").Value = -5336.1665

$ws.Range("H112").Value = 3858.4
$ws.Range("I112").Value = 1475.6666
$ws.Range("J112").Value = 4123.148
$ws.Range("K112").Value = 4426.9998
$ws.Range("L112").Value = 12369.444
$ws.Range("M112").Value = -3318.9998
$ws.Range("N112").Value = -14585.444

$ws.Range("H113").Value = 2916.26
$ws.Range("J113").Value = 3161.1777
$ws.Range("L113").Value = 9483.533100000001
$ws.Range("N113").Value = -13823.5331

$ws.Range("H131").Value = 868.11
$ws.Range("I131").Value = 687.5
$ws.Range("J131").Value = 875.63544
$ws.Range("K131").Value = 2062.5
$ws.Range("L131").Value = 2626.90632
$ws.Range("M131").Value = 2977.5
$ws.Range("N131").Value = -12706.90632

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4469.615
$ws.Range("I7").Value = 4540
$ws.Range("J7").Value = 4425.625
$ws.Range("K7").Value = 4540
$ws.Range("L7").Value = 4425.625
$ws.Range("M7").Value = -4428
$ws.Range("N7").Value = -4649.625

$ws.Range("H100").Value = 4537.778
$ws.Range("I100").Value = 3495
$ws.Range("J100").Value = 4835.7144
$ws.Range("K100").Value = 3495
$ws.Range("L100").Value = 4835.7144
$ws.Range("M100").Value = -2954
$ws.Range("N100").Value = -5917.7144

$ws.Range("H122").Value = 11647953
$ws.Range("I122").Value = 13976093
$ws.Range("J122").Value = 7251.25
$ws.Range("K122").Value = 41928279
$ws.Range("L122").Value = 21753.75
$ws.Range("M122").Value = -41925829
$ws.Range("N122").Value = -26653.75

$ws.Range("H126").Value = 4469.615
$ws.Range("I126").Value = 4540
$ws.Range("J126").Value = 4425.625
$ws.Range("K126").Value = 13620
$ws.Range("L126").Value = 13276.875
$ws.Range("M126").Value = -11150
$ws.Range("N126").Value = -18216.875

$ws.Range("H132").Value = 3687.5789
$ws.Range("I132").Value = 3337.2666
$ws.Range("J132").Value = 5001.25
$ws.Range("K132").Value = 10011.7998
$ws.Range("L132").Value = 15003.75
$ws.Range("M132").Value = -7481.799800000001
$ws.Range("N132").Value = -20063.75

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

$ws.Range("H139").Value = 66880
$ws.Range("J139").Value = 66880
$ws.Range("L139").Value = 66880
$ws.Range("N139").Value = -77160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 55415.145
$ws.Range("J133").Value = 55415.145
$ws.Range("L133").Value = 55415.145
$ws.Range("N133").Value = -65535.145

$ws.Range("H138").Value = 120428.5
$ws.Range("J138").Value = 120428.5
$ws.Range("L138").Value = 120428.5
$ws.Range("N138").Value = -130708.5

$ws.Range("H141").Value = 94000
$ws.Range("J141").Value = 94000
$ws.Range("L141").Value = 94000
$ws.Range("N141").Value = -104360
